$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New logboek rows (24-42): date, hours (D), activity (G)
# Data taken from the commit diff.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 24; Date = 42358; Hours = 4;   Activity = "Programmeren" },
    @{ Row = 25; Date = 42359; Hours = 4;   Activity = "Programmeren" },
    @{ Row = 26; Date = 42360; Hours = 6;   Activity = "Programmeren" },
    @{ Row = 27; Date = 42361; Hours = 3;   Activity = "Programmeren" },
    @{ Row = 28; Date = 42362; Hours = 3;   Activity = "Programmeren" },
    @{ Row = 29; Date = 42363; Hours = 3;   Activity = "Programmeren" },
    @{ Row = 30; Date = 42365; Hours = 4;   Activity = "Programmeren" },
    @{ Row = 31; Date = 42366; Hours = 0.5; Activity = "Teambespreking"; NoHoursStyle = $true },
    @{ Row = 32; Date = 42366; Hours = 3;   Activity = "Programmeren" },
    @{ Row = 33; Date = 42367; Hours = 5;   Activity = "Programmeren" },
    @{ Row = 34; Date = 42368; Hours = 6;   Activity = "Programmeren" },
    @{ Row = 35; Date = 42369; Hours = 1;   Activity = "Teambespreking en programmeren" },
    @{ Row = 36; Date = 42372; Hours = 4;   Activity = "Programmeren + tutorials bekijken" },
    @{ Row = 37; Date = 42373; Hours = 6.5; Activity = "Groepsgesprek, programmeren + tutorials" },
    @{ Row = 38; Date = 42374; Hours = 6;   Activity = "Programmeren + nadenken" },
    @{ Row = 39; Date = 42375; Hours = 2;   Activity = "Begeleidergesprek + programmeren" },
    @{ Row = 40; Date = 42375; Hours = 2;   Activity = "Programmeren" },
    @{ Row = 41; Date = 42376; Hours = 3;   Activity = "Programmeren + groepsgesprek" },
    @{ Row = 42; Date = 42376; Hours = 1;   Activity = "Vraaggenereren voorbereiden en instructie schrijven" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Date (column A) - copy number format from a known date cell so it
    # reuses the existing date style instead of minting a new one.
    $ws.Range("A23").Copy() | Out-Null
    $ws.Range("A" + $rowNum).PasteSpecial(-4122) | Out-Null
    $ws.Range("A" + $rowNum).Value = $r.Date

    # Hours (column D) - always copy the style from a known source cell so
    # that the no-style row (31) doesn't propagate its style to the rows
    # that follow it.
    if ($r.NoHoursStyle) {
        $ws.Range("D13").Copy() | Out-Null
    } else {
        $ws.Range("D23").Copy() | Out-Null
    }
    $ws.Range("D" + $rowNum).PasteSpecial(-4122) | Out-Null
    $ws.Range("D" + $rowNum).Value = $r.Hours

    # Activity (column G)
    $ws.Range("G23").Copy() | Out-Null
    $ws.Range("G" + $rowNum).PasteSpecial(-4122) | Out-Null
    $ws.Range("G" + $rowNum).Value = $r.Activity
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Totals formula now sums the whole column instead of a fixed range.
# ---------------------------------------------------------------------------
$ws.Range("W2").Formula = "=SUM(D:D)"

# ---------------------------------------------------------------------------
# Selection moves to A43.
# ---------------------------------------------------------------------------
$ws.Range("A43").Select()

$excel.Calculate()
